# code refactoring and loan accounting and charges added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoan_Input")

# shortname: switch from the shared-string "kar4" to the literal numeric value 392
$ws.Range("B3").Value = 392

# nominalinterestratedefault: 12 -> 1
$ws.Range("B11").Value = 1

# New loan accounting / charge mapping rows (31-42)
$accountingRows = @(
    @{ Row = 31; Field = "fundsource";                    Value = "Cash" },
    @{ Row = 32; Field = "loanprotfolio";                 Value = "Loan portfolio " },
    @{ Row = 33; Field = "interestreceivable";            Value = "Interest Receivable " },
    @{ Row = 34; Field = "penaltiesreceivable";            Value = "Penalties Receivable " },
    @{ Row = 35; Field = "transferinsuspense";            Value = "Transfer in Suspence " },
    @{ Row = 36; Field = "feesreceivable";                Value = "Fees Receivable" },
    @{ Row = 37; Field = "incomefrominterest";            Value = "Income from interest" },
    @{ Row = 38; Field = "incomefrompenalties";           Value = "Income from penalties" },
    @{ Row = 39; Field = "incomefromfees";                Value = "Income from fees" },
    @{ Row = 40; Field = "incomefromrecoveryrepayments";  Value = "Income from recovery repayments" },
    @{ Row = 41; Field = "loseswrittenoff";                Value = "Losses Writtenoff " },
    @{ Row = 42; Field = "overpaymentliability";           Value = "Overpayment Liability" }
)

# Values first (column B), then field names (column A) -- matches the
# order the strings were originally authored in the shared string table.
foreach ($r in $accountingRows) {
    $ws.Range("B$($r.Row)").Value = $r.Value
}
foreach ($r in $accountingRows) {
    $ws.Range("A$($r.Row)").Value = $r.Field
}

foreach ($r in $accountingRows) {
    $rowNum = $r.Row
    # Reuse the existing label/value styles (row 10: A=label style, B=value style)
    $null = $ws.Range("A10").Copy()
    $null = $ws.Range("A$rowNum").PasteSpecial(-4122)
    $null = $ws.Range("B10").Copy()
    $null = $ws.Range("B$rowNum").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the view state to match the edited workbook
$null = $ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
